# Add 2022-Q3 data: insert a new "2022-Q3" sheet (with fresh fund-holding data)
# right before the existing "2022-Q2" sheet, and add a matching summary row at
# the top of the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal TEXT (not auto-coerced to a Number), while
# leaving the cell's style untouched (no lingering "@"/quote-prefix format).
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a 2022-Q3 row above 2022-Q2 / 2022-Q1 (which shift
#    down by one row each). Written back-to-front so we never overwrite data
#    we still need to copy forward.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 4 (was row 3 / 2022-Q1): clone A-column style, then old row-3 values.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 4
$total.Cells.Item(4, 4).Value = 0.54

# Row 3 (was row 2 / 2022-Q2): old row-2 values (A3 already holds style+1).
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 0.5

# Row 2: brand-new 2022-Q3 summary values (A2 keeps its existing style/value).
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.41

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, positioned right before "2022-Q2".
#    NOTE: once a sheet is added/renamed the old $q2-style object handles can
#    go stale (their index shifts under them), so re-look-up sheets by name
#    after any Add()/Name= operation before using them again.
# ---------------------------------------------------------------------------
$q2ref = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2ref)
$newSheet.Name = "2022-Q3"

$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Match page margins used by the other quarterly sheets (0.75/0.75/1/1/.5/.5 in).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Copy header row + column-A formatting from the 2022-Q2 sheet so the new
# sheet matches the existing look (bold header, bordered index column).
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$q2.Range("A2:A5").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

# Header labels (identical to the other quarterly sheets).
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Row 2: 014277 / 万家北交所慧选两年定期开放混合A
$q3.Cells.Item(2, 1).Value = 0
Set-TextValue $q3.Cells.Item(2, 2) "014277"
$q3.Cells.Item(2, 3).Value = "万家北交所慧选两年定期开放混合A"
Set-TextValue $q3.Cells.Item(2, 4) "3.56"
Set-TextValue $q3.Cells.Item(2, 5) "93.97"
Set-TextValue $q3.Cells.Item(2, 6) "6.05"
Set-TextValue $q3.Cells.Item(2, 7) "0.2154"
$q3.Cells.Item(2, 8).Value = 5

# Row 3: 014279 / 汇添富北交所创新精选两年定开混合A
$q3.Cells.Item(3, 1).Value = 1
Set-TextValue $q3.Cells.Item(3, 2) "014279"
$q3.Cells.Item(3, 3).Value = "汇添富北交所创新精选两年定开混合A"
Set-TextValue $q3.Cells.Item(3, 4) "3.20"
Set-TextValue $q3.Cells.Item(3, 5) "93.27"
Set-TextValue $q3.Cells.Item(3, 6) "4.34"
Set-TextValue $q3.Cells.Item(3, 7) "0.1389"
$q3.Cells.Item(3, 8).Value = 8

# Row 4: 014278 / 万家北交所慧选两年定期开放混合C
$q3.Cells.Item(4, 1).Value = 2
Set-TextValue $q3.Cells.Item(4, 2) "014278"
$q3.Cells.Item(4, 3).Value = "万家北交所慧选两年定期开放混合C"
Set-TextValue $q3.Cells.Item(4, 4) "0.49"
Set-TextValue $q3.Cells.Item(4, 5) "93.97"
Set-TextValue $q3.Cells.Item(4, 6) "6.05"
Set-TextValue $q3.Cells.Item(4, 7) "0.0296"
$q3.Cells.Item(4, 8).Value = 5

# Row 5: 014280 / 汇添富北交所创新精选两年定开混合C
$q3.Cells.Item(5, 1).Value = 3
Set-TextValue $q3.Cells.Item(5, 2) "014280"
$q3.Cells.Item(5, 3).Value = "汇添富北交所创新精选两年定开混合C"
Set-TextValue $q3.Cells.Item(5, 4) "0.51"
Set-TextValue $q3.Cells.Item(5, 5) "93.27"
Set-TextValue $q3.Cells.Item(5, 6) "4.34"
Set-TextValue $q3.Cells.Item(5, 7) "0.0221"
$q3.Cells.Item(5, 8).Value = 8

# ---------------------------------------------------------------------------
# 3) Restore original tab-selection state: the "2022-Q1" sheet (now the 4th
#    tab) was the selected tab before this edit.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Select()
